# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row number (key) -> new F-column value (value)
$updates = @{
    4  = 103
    7  = 59
    8  = 489
    9  = 6423
    10 = 171
    11 = 131
    12 = 1017
    13 = 320
    14 = 102
    15 = 177
    16 = 467
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
